# Update rows 556-677: shift the weekly price records down by two rows
# (row R now holds the data that used to live at row R-2), and give the
# first two rows (556-557) a brand-new price record.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(556, 4).Value = 44798
$ws.Cells.Item(556, 14).Value = 18000
$ws.Cells.Item(556, 15).Value = 18000
$ws.Cells.Item(556, 16).Value = 18000
$ws.Cells.Item(556, 19).Value = 900
$ws.Cells.Item(557, 4).Value = 44798
$ws.Cells.Item(557, 13).Value = 160
$ws.Cells.Item(557, 14).Value = 19000
$ws.Cells.Item(557, 15).Value = 20000
$ws.Cells.Item(557, 16).Value = 19500
$ws.Cells.Item(557, 19).Value = 975
$ws.Cells.Item(558, 4).Value = 44342
$ws.Cells.Item(558, 13).Value = 80
$ws.Cells.Item(558, 14).Value = 11000
$ws.Cells.Item(558, 15).Value = 11000
$ws.Cells.Item(558, 16).Value = 11000
$ws.Cells.Item(558, 19).Value = 550
$ws.Cells.Item(559, 4).Value = 44342
$ws.Cells.Item(559, 14).Value = 12000
$ws.Cells.Item(559, 15).Value = 13000
$ws.Cells.Item(559, 16).Value = 12500
$ws.Cells.Item(559, 19).Value = 625
$ws.Cells.Item(560, 4).Value = 44678
$ws.Cells.Item(560, 13).Value = 160
$ws.Cells.Item(560, 14).Value = 12000
$ws.Cells.Item(560, 15).Value = 13000
$ws.Cells.Item(560, 16).Value = 12500
$ws.Cells.Item(560, 19).Value = 625
$ws.Cells.Item(561, 4).Value = 44678
$ws.Cells.Item(561, 14).Value = 14000
$ws.Cells.Item(561, 15).Value = 15000
$ws.Cells.Item(561, 16).Value = 14500
$ws.Cells.Item(561, 19).Value = 725
$ws.Cells.Item(562, 4).Value = 44551
$ws.Cells.Item(562, 14).Value = 11000
$ws.Cells.Item(562, 15).Value = 11000
$ws.Cells.Item(562, 16).Value = 11000
$ws.Cells.Item(562, 19).Value = 550
$ws.Cells.Item(563, 4).Value = 44551
$ws.Cells.Item(563, 13).Value = 240
$ws.Cells.Item(563, 14).Value = 12000
$ws.Cells.Item(563, 15).Value = 13000
$ws.Cells.Item(563, 16).Value = 12500
$ws.Cells.Item(563, 19).Value = 625
$ws.Cells.Item(564, 4).Value = 44648
$ws.Cells.Item(564, 14).Value = 18000
$ws.Cells.Item(564, 15).Value = 18000
$ws.Cells.Item(564, 16).Value = 18000
$ws.Cells.Item(564, 19).Value = 900
$ws.Cells.Item(565, 4).Value = 44648
$ws.Cells.Item(565, 13).Value = 160
$ws.Cells.Item(565, 14).Value = 19000
$ws.Cells.Item(565, 15).Value = 20000
$ws.Cells.Item(565, 16).Value = 19500
$ws.Cells.Item(565, 19).Value = 975
$ws.Cells.Item(566, 4).Value = 44291
$ws.Cells.Item(566, 13).Value = 80
$ws.Cells.Item(567, 4).Value = 44291
$ws.Cells.Item(567, 13).Value = 240
$ws.Cells.Item(568, 4).Value = 44585
$ws.Cells.Item(568, 13).Value = 100
$ws.Cells.Item(568, 14).Value = 13000
$ws.Cells.Item(568, 15).Value = 13000
$ws.Cells.Item(568, 16).Value = 13000
$ws.Cells.Item(568, 19).Value = 650
$ws.Cells.Item(569, 4).Value = 44585
$ws.Cells.Item(569, 13).Value = 200
$ws.Cells.Item(569, 14).Value = 14000
$ws.Cells.Item(569, 15).Value = 15000
$ws.Cells.Item(569, 16).Value = 14500
$ws.Cells.Item(569, 19).Value = 725
$ws.Cells.Item(570, 4).Value = 44376
$ws.Cells.Item(570, 14).Value = 11000
$ws.Cells.Item(570, 15).Value = 11000
$ws.Cells.Item(570, 16).Value = 11000
$ws.Cells.Item(570, 19).Value = 550
$ws.Cells.Item(571, 4).Value = 44376
$ws.Cells.Item(571, 13).Value = 240
$ws.Cells.Item(571, 14).Value = 12000
$ws.Cells.Item(571, 15).Value = 13000
$ws.Cells.Item(571, 16).Value = 12500
$ws.Cells.Item(571, 19).Value = 625
$ws.Cells.Item(572, 4).Value = 44242
$ws.Cells.Item(572, 13).Value = 80
$ws.Cells.Item(572, 14).Value = 8000
$ws.Cells.Item(572, 15).Value = 8000
$ws.Cells.Item(572, 16).Value = 8000
$ws.Cells.Item(572, 19).Value = 400
$ws.Cells.Item(573, 4).Value = 44242
$ws.Cells.Item(573, 13).Value = 300
$ws.Cells.Item(573, 14).Value = 9000
$ws.Cells.Item(573, 15).Value = 10000
$ws.Cells.Item(573, 16).Value = 9500
$ws.Cells.Item(573, 19).Value = 475
$ws.Cells.Item(574, 4).Value = 44391
$ws.Cells.Item(574, 13).Value = 120
$ws.Cells.Item(574, 14).Value = 11000
$ws.Cells.Item(574, 15).Value = 11000
$ws.Cells.Item(574, 16).Value = 11000
$ws.Cells.Item(574, 19).Value = 550
$ws.Cells.Item(575, 4).Value = 44391
$ws.Cells.Item(575, 13).Value = 240
$ws.Cells.Item(575, 14).Value = 12000
$ws.Cells.Item(575, 16).Value = 12500
$ws.Cells.Item(575, 19).Value = 625
$ws.Cells.Item(576, 4).Value = 44600
$ws.Cells.Item(576, 13).Value = 150
$ws.Cells.Item(576, 14).Value = 12000
$ws.Cells.Item(576, 15).Value = 12000
$ws.Cells.Item(576, 16).Value = 12000
$ws.Cells.Item(576, 19).Value = 600
$ws.Cells.Item(577, 4).Value = 44600
$ws.Cells.Item(577, 12).Value = "Primera Pintón"
$ws.Cells.Item(577, 13).Value = 500
$ws.Cells.Item(577, 14).Value = 12500
$ws.Cells.Item(577, 15).Value = 13000
$ws.Cells.Item(577, 16).Value = 12750
$ws.Cells.Item(577, 19).Value = 638
$ws.Cells.Item(578, 4).Value = 44763
$ws.Cells.Item(578, 12).Value = "Pintón"
$ws.Cells.Item(578, 13).Value = 80
$ws.Cells.Item(578, 14).Value = 29000
$ws.Cells.Item(578, 15).Value = 29000
$ws.Cells.Item(578, 16).Value = 29000
$ws.Cells.Item(578, 19).Value = 1450
$ws.Cells.Item(579, 4).Value = 44371
$ws.Cells.Item(579, 13).Value = 120
$ws.Cells.Item(580, 4).Value = 44371
$ws.Cells.Item(580, 13).Value = 240
$ws.Cells.Item(581, 4).Value = 44355
$ws.Cells.Item(581, 14).Value = 11000
$ws.Cells.Item(581, 15).Value = 11000
$ws.Cells.Item(581, 16).Value = 11000
$ws.Cells.Item(581, 19).Value = 550
$ws.Cells.Item(582, 4).Value = 44355
$ws.Cells.Item(582, 13).Value = 160
$ws.Cells.Item(582, 14).Value = 12000
$ws.Cells.Item(582, 15).Value = 13000
$ws.Cells.Item(582, 16).Value = 12500
$ws.Cells.Item(582, 19).Value = 625
$ws.Cells.Item(583, 4).Value = 44579
$ws.Cells.Item(583, 12).Value = "Pintón"
$ws.Cells.Item(583, 13).Value = 80
$ws.Cells.Item(583, 14).Value = 14000
$ws.Cells.Item(583, 15).Value = 14000
$ws.Cells.Item(583, 16).Value = 14000
$ws.Cells.Item(583, 19).Value = 700
$ws.Cells.Item(584, 4).Value = 44579
$ws.Cells.Item(584, 12).Value = "Primera Pintón"
$ws.Cells.Item(584, 13).Value = 240
$ws.Cells.Item(584, 14).Value = 15000
$ws.Cells.Item(584, 15).Value = 16000
$ws.Cells.Item(584, 16).Value = 15500
$ws.Cells.Item(584, 19).Value = 775
$ws.Cells.Item(585, 4).Value = 44685
$ws.Cells.Item(585, 14).Value = 17000
$ws.Cells.Item(585, 15).Value = 18000
$ws.Cells.Item(585, 16).Value = 17500
$ws.Cells.Item(585, 19).Value = 875
$ws.Cells.Item(586, 4).Value = 44558
$ws.Cells.Item(586, 13).Value = 80
$ws.Cells.Item(586, 14).Value = 12000
$ws.Cells.Item(586, 15).Value = 12000
$ws.Cells.Item(586, 16).Value = 12000
$ws.Cells.Item(586, 19).Value = 600
$ws.Cells.Item(587, 4).Value = 44558
$ws.Cells.Item(587, 13).Value = 240
$ws.Cells.Item(587, 14).Value = 13000
$ws.Cells.Item(587, 15).Value = 14000
$ws.Cells.Item(587, 16).Value = 13500
$ws.Cells.Item(587, 19).Value = 675
$ws.Cells.Item(588, 4).Value = 44434
$ws.Cells.Item(588, 12).Value = "Pintón"
$ws.Cells.Item(588, 13).Value = 300
$ws.Cells.Item(588, 14).Value = 14500
$ws.Cells.Item(588, 15).Value = 15000
$ws.Cells.Item(588, 16).Value = 14750
$ws.Cells.Item(588, 19).Value = 738
$ws.Cells.Item(589, 4).Value = 44434
$ws.Cells.Item(589, 12).Value = "Primera Pintón"
$ws.Cells.Item(589, 13).Value = 600
$ws.Cells.Item(589, 14).Value = 15500
$ws.Cells.Item(589, 15).Value = 16000
$ws.Cells.Item(589, 16).Value = 15750
$ws.Cells.Item(589, 19).Value = 788
$ws.Cells.Item(590, 4).Value = 44278
$ws.Cells.Item(590, 13).Value = 360
$ws.Cells.Item(590, 14).Value = 12000
$ws.Cells.Item(590, 15).Value = 13000
$ws.Cells.Item(590, 16).Value = 12500
$ws.Cells.Item(590, 19).Value = 625
$ws.Cells.Item(591, 4).Value = 44442
$ws.Cells.Item(591, 13).Value = 300
$ws.Cells.Item(591, 14).Value = 13500
$ws.Cells.Item(591, 16).Value = 13750
$ws.Cells.Item(591, 19).Value = 688
$ws.Cells.Item(592, 4).Value = 44442
$ws.Cells.Item(592, 13).Value = 300
$ws.Cells.Item(592, 14).Value = 14500
$ws.Cells.Item(592, 15).Value = 15000
$ws.Cells.Item(592, 16).Value = 14750
$ws.Cells.Item(592, 19).Value = 738
$ws.Cells.Item(593, 4).Value = 44238
$ws.Cells.Item(593, 13).Value = 120
$ws.Cells.Item(593, 14).Value = 14000
$ws.Cells.Item(593, 15).Value = 14000
$ws.Cells.Item(593, 16).Value = 14000
$ws.Cells.Item(593, 19).Value = 700
$ws.Cells.Item(594, 4).Value = 44238
$ws.Cells.Item(594, 13).Value = 240
$ws.Cells.Item(594, 14).Value = 15000
$ws.Cells.Item(594, 15).Value = 16000
$ws.Cells.Item(594, 16).Value = 15500
$ws.Cells.Item(594, 19).Value = 775
$ws.Cells.Item(595, 4).Value = 44781
$ws.Cells.Item(595, 12).Value = "Pintón"
$ws.Cells.Item(595, 13).Value = 80
$ws.Cells.Item(595, 14).Value = 21000
$ws.Cells.Item(595, 15).Value = 21000
$ws.Cells.Item(595, 16).Value = 21000
$ws.Cells.Item(595, 19).Value = 1050
$ws.Cells.Item(596, 4).Value = 44781
$ws.Cells.Item(596, 12).Value = "Primera Pintón"
$ws.Cells.Item(596, 13).Value = 160
$ws.Cells.Item(596, 14).Value = 22000
$ws.Cells.Item(596, 15).Value = 23000
$ws.Cells.Item(596, 16).Value = 22500
$ws.Cells.Item(596, 19).Value = 1125
$ws.Cells.Item(597, 4).Value = 44336
$ws.Cells.Item(597, 13).Value = 300
$ws.Cells.Item(598, 4).Value = 44343
$ws.Cells.Item(598, 13).Value = 80
$ws.Cells.Item(598, 14).Value = 11000
$ws.Cells.Item(598, 15).Value = 11000
$ws.Cells.Item(598, 16).Value = 11000
$ws.Cells.Item(598, 19).Value = 550
$ws.Cells.Item(599, 4).Value = 44343
$ws.Cells.Item(599, 13).Value = 160
$ws.Cells.Item(599, 14).Value = 12000
$ws.Cells.Item(599, 15).Value = 13000
$ws.Cells.Item(599, 16).Value = 12500
$ws.Cells.Item(599, 19).Value = 625
$ws.Cells.Item(600, 4).Value = 44533
$ws.Cells.Item(600, 14).Value = 18000
$ws.Cells.Item(600, 15).Value = 18000
$ws.Cells.Item(600, 16).Value = 18000
$ws.Cells.Item(600, 19).Value = 900
$ws.Cells.Item(601, 4).Value = 44533
$ws.Cells.Item(601, 12).Value = "Primera Pintón"
$ws.Cells.Item(601, 14).Value = 19000
$ws.Cells.Item(601, 15).Value = 20000
$ws.Cells.Item(601, 16).Value = 19500
$ws.Cells.Item(601, 19).Value = 975
$ws.Cells.Item(602, 4).Value = 44365
$ws.Cells.Item(602, 12).Value = "Pintón"
$ws.Cells.Item(602, 13).Value = 120
$ws.Cells.Item(602, 14).Value = 14500
$ws.Cells.Item(602, 15).Value = 15000
$ws.Cells.Item(602, 16).Value = 14750
$ws.Cells.Item(602, 19).Value = 738
$ws.Cells.Item(603, 4).Value = 44454
$ws.Cells.Item(603, 13).Value = 240
$ws.Cells.Item(603, 14).Value = 18000
$ws.Cells.Item(603, 15).Value = 19000
$ws.Cells.Item(603, 16).Value = 18500
$ws.Cells.Item(603, 19).Value = 925
$ws.Cells.Item(604, 4).Value = 44454
$ws.Cells.Item(604, 13).Value = 300
$ws.Cells.Item(604, 14).Value = 20000
$ws.Cells.Item(604, 15).Value = 21000
$ws.Cells.Item(604, 16).Value = 20500
$ws.Cells.Item(604, 19).Value = 1025
$ws.Cells.Item(605, 4).Value = 44561
$ws.Cells.Item(605, 13).Value = 80
$ws.Cells.Item(605, 14).Value = 11000
$ws.Cells.Item(605, 15).Value = 11000
$ws.Cells.Item(605, 16).Value = 11000
$ws.Cells.Item(605, 19).Value = 550
$ws.Cells.Item(606, 4).Value = 44561
$ws.Cells.Item(606, 13).Value = 160
$ws.Cells.Item(606, 14).Value = 12000
$ws.Cells.Item(606, 15).Value = 13000
$ws.Cells.Item(606, 16).Value = 12500
$ws.Cells.Item(606, 19).Value = 625
$ws.Cells.Item(607, 4).Value = 44421
$ws.Cells.Item(607, 14).Value = 12000
$ws.Cells.Item(607, 15).Value = 12000
$ws.Cells.Item(607, 16).Value = 12000
$ws.Cells.Item(607, 19).Value = 600
$ws.Cells.Item(608, 4).Value = 44421
$ws.Cells.Item(608, 14).Value = 13000
$ws.Cells.Item(608, 15).Value = 14000
$ws.Cells.Item(608, 16).Value = 13500
$ws.Cells.Item(608, 19).Value = 675
$ws.Cells.Item(609, 4).Value = 44419
$ws.Cells.Item(609, 13).Value = 300
$ws.Cells.Item(609, 14).Value = 11000
$ws.Cells.Item(609, 15).Value = 11000
$ws.Cells.Item(609, 16).Value = 11000
$ws.Cells.Item(609, 19).Value = 550
$ws.Cells.Item(610, 4).Value = 44419
$ws.Cells.Item(610, 13).Value = 600
$ws.Cells.Item(610, 14).Value = 12000
$ws.Cells.Item(610, 15).Value = 13000
$ws.Cells.Item(610, 16).Value = 12500
$ws.Cells.Item(610, 19).Value = 625
$ws.Cells.Item(611, 4).Value = 44669
$ws.Cells.Item(611, 12).Value = "Pintón"
$ws.Cells.Item(611, 13).Value = 80
$ws.Cells.Item(611, 14).Value = 15000
$ws.Cells.Item(611, 15).Value = 15000
$ws.Cells.Item(611, 16).Value = 15000
$ws.Cells.Item(611, 19).Value = 750
$ws.Cells.Item(612, 4).Value = 44669
$ws.Cells.Item(612, 12).Value = "Primera Pintón"
$ws.Cells.Item(612, 13).Value = 160
$ws.Cells.Item(612, 14).Value = 16000
$ws.Cells.Item(612, 15).Value = 17000
$ws.Cells.Item(612, 16).Value = 16500
$ws.Cells.Item(612, 19).Value = 825
$ws.Cells.Item(613, 4).Value = 44489
$ws.Cells.Item(613, 13).Value = 400
$ws.Cells.Item(613, 14).Value = 26000
$ws.Cells.Item(613, 15).Value = 27000
$ws.Cells.Item(613, 16).Value = 26500
$ws.Cells.Item(613, 19).Value = 1325
$ws.Cells.Item(614, 4).Value = 44699
$ws.Cells.Item(614, 13).Value = 80
$ws.Cells.Item(614, 14).Value = 11000
$ws.Cells.Item(614, 15).Value = 11000
$ws.Cells.Item(614, 16).Value = 11000
$ws.Cells.Item(614, 19).Value = 550
$ws.Cells.Item(615, 4).Value = 44699
$ws.Cells.Item(615, 13).Value = 160
$ws.Cells.Item(615, 14).Value = 12000
$ws.Cells.Item(615, 15).Value = 13000
$ws.Cells.Item(615, 16).Value = 12500
$ws.Cells.Item(615, 19).Value = 625
$ws.Cells.Item(616, 4).Value = 44636
$ws.Cells.Item(616, 12).Value = "Pintón"
$ws.Cells.Item(616, 13).Value = 100
$ws.Cells.Item(616, 15).Value = 17000
$ws.Cells.Item(616, 16).Value = 17000
$ws.Cells.Item(616, 19).Value = 850
$ws.Cells.Item(617, 4).Value = 44636
$ws.Cells.Item(617, 12).Value = "Primera Pintón"
$ws.Cells.Item(617, 13).Value = 240
$ws.Cells.Item(617, 14).Value = 18000
$ws.Cells.Item(617, 15).Value = 19000
$ws.Cells.Item(617, 16).Value = 18500
$ws.Cells.Item(617, 19).Value = 925
$ws.Cells.Item(618, 4).Value = 44405
$ws.Cells.Item(618, 12).Value = "Primera Pintón"
$ws.Cells.Item(618, 13).Value = 240
$ws.Cells.Item(618, 14).Value = 17000
$ws.Cells.Item(618, 15).Value = 18000
$ws.Cells.Item(618, 16).Value = 17500
$ws.Cells.Item(618, 19).Value = 875
$ws.Cells.Item(619, 4).Value = 44202
$ws.Cells.Item(619, 12).Value = "Pintón"
$ws.Cells.Item(619, 13).Value = 300
$ws.Cells.Item(619, 14).Value = 15000
$ws.Cells.Item(619, 15).Value = 16000
$ws.Cells.Item(619, 16).Value = 15500
$ws.Cells.Item(619, 19).Value = 775
$ws.Cells.Item(620, 4).Value = 44273
$ws.Cells.Item(620, 13).Value = 120
$ws.Cells.Item(620, 14).Value = 12000
$ws.Cells.Item(620, 15).Value = 12000
$ws.Cells.Item(620, 16).Value = 12000
$ws.Cells.Item(620, 19).Value = 600
$ws.Cells.Item(621, 4).Value = 44273
$ws.Cells.Item(621, 13).Value = 240
$ws.Cells.Item(621, 14).Value = 13000
$ws.Cells.Item(621, 15).Value = 14000
$ws.Cells.Item(621, 16).Value = 13500
$ws.Cells.Item(621, 19).Value = 675
$ws.Cells.Item(622, 4).Value = 44777
$ws.Cells.Item(622, 14).Value = 25000
$ws.Cells.Item(622, 15).Value = 25000
$ws.Cells.Item(622, 16).Value = 25000
$ws.Cells.Item(622, 19).Value = 1250
$ws.Cells.Item(623, 4).Value = 44777
$ws.Cells.Item(623, 14).Value = 26000
$ws.Cells.Item(623, 15).Value = 27000
$ws.Cells.Item(623, 16).Value = 26500
$ws.Cells.Item(623, 19).Value = 1325
$ws.Cells.Item(624, 4).Value = 44159
$ws.Cells.Item(624, 13).Value = 80
$ws.Cells.Item(624, 14).Value = 17000
$ws.Cells.Item(624, 15).Value = 17000
$ws.Cells.Item(624, 16).Value = 17000
$ws.Cells.Item(624, 19).Value = 850
$ws.Cells.Item(625, 4).Value = 44159
$ws.Cells.Item(625, 13).Value = 160
$ws.Cells.Item(625, 14).Value = 18000
$ws.Cells.Item(625, 15).Value = 19000
$ws.Cells.Item(625, 16).Value = 18500
$ws.Cells.Item(625, 19).Value = 925
$ws.Cells.Item(626, 4).Value = 44589
$ws.Cells.Item(626, 13).Value = 400
$ws.Cells.Item(626, 14).Value = 11000
$ws.Cells.Item(626, 15).Value = 11000
$ws.Cells.Item(626, 16).Value = 11000
$ws.Cells.Item(626, 19).Value = 550
$ws.Cells.Item(627, 4).Value = 44589
$ws.Cells.Item(627, 13).Value = 400
$ws.Cells.Item(627, 14).Value = 12000
$ws.Cells.Item(627, 15).Value = 13000
$ws.Cells.Item(627, 16).Value = 12500
$ws.Cells.Item(627, 19).Value = 625
$ws.Cells.Item(628, 4).Value = 44263
$ws.Cells.Item(628, 12).Value = "Pintón"
$ws.Cells.Item(628, 13).Value = 80
$ws.Cells.Item(628, 14).Value = 14000
$ws.Cells.Item(628, 16).Value = 14000
$ws.Cells.Item(628, 19).Value = 700
$ws.Cells.Item(629, 4).Value = 44263
$ws.Cells.Item(629, 12).Value = "Primera Pintón"
$ws.Cells.Item(629, 13).Value = 160
$ws.Cells.Item(629, 14).Value = 15000
$ws.Cells.Item(629, 15).Value = 16000
$ws.Cells.Item(629, 16).Value = 15500
$ws.Cells.Item(629, 19).Value = 775
$ws.Cells.Item(630, 4).Value = 44309
$ws.Cells.Item(630, 14).Value = 13000
$ws.Cells.Item(630, 15).Value = 14000
$ws.Cells.Item(630, 16).Value = 13500
$ws.Cells.Item(630, 19).Value = 675
$ws.Cells.Item(631, 4).Value = 44771
$ws.Cells.Item(631, 13).Value = 80
$ws.Cells.Item(631, 14).Value = 30000
$ws.Cells.Item(631, 15).Value = 30000
$ws.Cells.Item(631, 16).Value = 30000
$ws.Cells.Item(631, 19).Value = 1500
$ws.Cells.Item(632, 4).Value = 44267
$ws.Cells.Item(632, 13).Value = 360
$ws.Cells.Item(632, 14).Value = 14000
$ws.Cells.Item(632, 15).Value = 15000
$ws.Cells.Item(632, 16).Value = 14500
$ws.Cells.Item(632, 19).Value = 725
$ws.Cells.Item(633, 4).Value = 44413
$ws.Cells.Item(633, 13).Value = 200
$ws.Cells.Item(633, 14).Value = 11000
$ws.Cells.Item(633, 15).Value = 11000
$ws.Cells.Item(633, 16).Value = 11000
$ws.Cells.Item(633, 19).Value = 550
$ws.Cells.Item(634, 4).Value = 44413
$ws.Cells.Item(634, 13).Value = 600
$ws.Cells.Item(634, 14).Value = 12000
$ws.Cells.Item(634, 15).Value = 13000
$ws.Cells.Item(634, 16).Value = 12500
$ws.Cells.Item(634, 19).Value = 625
$ws.Cells.Item(635, 4).Value = 44515
$ws.Cells.Item(635, 13).Value = 160
$ws.Cells.Item(635, 14).Value = 15500
$ws.Cells.Item(635, 15).Value = 16000
$ws.Cells.Item(635, 16).Value = 15750
$ws.Cells.Item(635, 19).Value = 788
$ws.Cells.Item(636, 4).Value = 44515
$ws.Cells.Item(636, 13).Value = 240
$ws.Cells.Item(636, 14).Value = 17000
$ws.Cells.Item(636, 15).Value = 18000
$ws.Cells.Item(636, 16).Value = 17500
$ws.Cells.Item(636, 19).Value = 875
$ws.Cells.Item(637, 4).Value = 44356
$ws.Cells.Item(637, 13).Value = 60
$ws.Cells.Item(637, 14).Value = 11000
$ws.Cells.Item(637, 15).Value = 11000
$ws.Cells.Item(637, 16).Value = 11000
$ws.Cells.Item(637, 19).Value = 550
$ws.Cells.Item(638, 4).Value = 44356
$ws.Cells.Item(638, 13).Value = 120
$ws.Cells.Item(638, 14).Value = 12000
$ws.Cells.Item(638, 15).Value = 13000
$ws.Cells.Item(638, 16).Value = 12500
$ws.Cells.Item(638, 19).Value = 625
$ws.Cells.Item(639, 4).Value = 44379
$ws.Cells.Item(639, 13).Value = 120
$ws.Cells.Item(639, 14).Value = 10000
$ws.Cells.Item(639, 15).Value = 10000
$ws.Cells.Item(639, 16).Value = 10000
$ws.Cells.Item(639, 19).Value = 500
$ws.Cells.Item(640, 4).Value = 44379
$ws.Cells.Item(640, 12).Value = "Primera Pintón"
$ws.Cells.Item(640, 13).Value = 240
$ws.Cells.Item(640, 14).Value = 11000
$ws.Cells.Item(640, 16).Value = 11500
$ws.Cells.Item(640, 19).Value = 575
$ws.Cells.Item(641, 4).Value = 44322
$ws.Cells.Item(641, 12).Value = "Pintón"
$ws.Cells.Item(641, 14).Value = 14000
$ws.Cells.Item(641, 15).Value = 15000
$ws.Cells.Item(641, 16).Value = 14500
$ws.Cells.Item(641, 19).Value = 725
$ws.Cells.Item(642, 4).Value = 44221
$ws.Cells.Item(642, 13).Value = 80
$ws.Cells.Item(642, 14).Value = 12000
$ws.Cells.Item(642, 15).Value = 12000
$ws.Cells.Item(642, 16).Value = 12000
$ws.Cells.Item(642, 19).Value = 600
$ws.Cells.Item(643, 4).Value = 44221
$ws.Cells.Item(643, 13).Value = 240
$ws.Cells.Item(643, 14).Value = 13000
$ws.Cells.Item(643, 15).Value = 14000
$ws.Cells.Item(643, 16).Value = 13500
$ws.Cells.Item(643, 19).Value = 675
$ws.Cells.Item(644, 4).Value = 44497
$ws.Cells.Item(644, 13).Value = 240
$ws.Cells.Item(644, 14).Value = 19000
$ws.Cells.Item(644, 15).Value = 20000
$ws.Cells.Item(644, 16).Value = 19500
$ws.Cells.Item(644, 19).Value = 975
$ws.Cells.Item(645, 4).Value = 44497
$ws.Cells.Item(645, 13).Value = 160
$ws.Cells.Item(645, 14).Value = 21000
$ws.Cells.Item(645, 15).Value = 22000
$ws.Cells.Item(645, 16).Value = 21500
$ws.Cells.Item(645, 19).Value = 1075
$ws.Cells.Item(646, 4).Value = 44782
$ws.Cells.Item(646, 13).Value = 60
$ws.Cells.Item(646, 14).Value = 21000
$ws.Cells.Item(646, 15).Value = 21000
$ws.Cells.Item(646, 16).Value = 21000
$ws.Cells.Item(646, 19).Value = 1050
$ws.Cells.Item(647, 4).Value = 44782
$ws.Cells.Item(647, 13).Value = 120
$ws.Cells.Item(647, 14).Value = 22000
$ws.Cells.Item(647, 15).Value = 23000
$ws.Cells.Item(647, 16).Value = 22500
$ws.Cells.Item(647, 19).Value = 1125
$ws.Cells.Item(648, 4).Value = 44435
$ws.Cells.Item(648, 13).Value = 1340
$ws.Cells.Item(648, 14).Value = 11500
$ws.Cells.Item(648, 15).Value = 15000
$ws.Cells.Item(648, 16).Value = 13177
$ws.Cells.Item(648, 19).Value = 659
$ws.Cells.Item(649, 4).Value = 44435
$ws.Cells.Item(649, 13).Value = 2020
$ws.Cells.Item(649, 14).Value = 12500
$ws.Cells.Item(649, 15).Value = 16000
$ws.Cells.Item(649, 16).Value = 14505
$ws.Cells.Item(649, 19).Value = 725
$ws.Cells.Item(650, 4).Value = 44251
$ws.Cells.Item(650, 13).Value = 660
$ws.Cells.Item(650, 14).Value = 9000
$ws.Cells.Item(650, 15).Value = 9500
$ws.Cells.Item(650, 16).Value = 9273
$ws.Cells.Item(650, 19).Value = 464
$ws.Cells.Item(651, 4).Value = 44251
$ws.Cells.Item(651, 13).Value = 300
$ws.Cells.Item(651, 14).Value = 10000
$ws.Cells.Item(651, 15).Value = 11000
$ws.Cells.Item(651, 16).Value = 10667
$ws.Cells.Item(651, 19).Value = 533
$ws.Cells.Item(652, 4).Value = 44319
$ws.Cells.Item(652, 14).Value = 13000
$ws.Cells.Item(652, 15).Value = 13000
$ws.Cells.Item(652, 16).Value = 13000
$ws.Cells.Item(652, 19).Value = 650
$ws.Cells.Item(653, 4).Value = 44319
$ws.Cells.Item(653, 13).Value = 240
$ws.Cells.Item(653, 14).Value = 14500
$ws.Cells.Item(653, 15).Value = 15000
$ws.Cells.Item(653, 16).Value = 14750
$ws.Cells.Item(653, 19).Value = 738
$ws.Cells.Item(654, 4).Value = 44344
$ws.Cells.Item(654, 13).Value = 80
$ws.Cells.Item(654, 14).Value = 12000
$ws.Cells.Item(654, 15).Value = 12000
$ws.Cells.Item(654, 16).Value = 12000
$ws.Cells.Item(654, 19).Value = 600
$ws.Cells.Item(655, 4).Value = 44344
$ws.Cells.Item(655, 13).Value = 160
$ws.Cells.Item(655, 14).Value = 13000
$ws.Cells.Item(655, 15).Value = 14000
$ws.Cells.Item(655, 16).Value = 13500
$ws.Cells.Item(655, 19).Value = 675
$ws.Cells.Item(656, 4).Value = 44232
$ws.Cells.Item(656, 13).Value = 760
$ws.Cells.Item(656, 14).Value = 14500
$ws.Cells.Item(656, 15).Value = 15000
$ws.Cells.Item(656, 16).Value = 14789
$ws.Cells.Item(656, 19).Value = 739
$ws.Cells.Item(657, 4).Value = 44232
$ws.Cells.Item(657, 14).Value = 15500
$ws.Cells.Item(657, 15).Value = 16000
$ws.Cells.Item(657, 16).Value = 15750
$ws.Cells.Item(657, 19).Value = 788
$ws.Cells.Item(658, 4).Value = 44455
$ws.Cells.Item(658, 14).Value = 18000
$ws.Cells.Item(658, 15).Value = 19000
$ws.Cells.Item(658, 16).Value = 18500
$ws.Cells.Item(658, 19).Value = 925
$ws.Cells.Item(659, 4).Value = 44455
$ws.Cells.Item(659, 14).Value = 20000
$ws.Cells.Item(659, 15).Value = 21000
$ws.Cells.Item(659, 16).Value = 20500
$ws.Cells.Item(659, 19).Value = 1025
$ws.Cells.Item(660, 4).Value = 44484
$ws.Cells.Item(660, 13).Value = 300
$ws.Cells.Item(660, 14).Value = 19500
$ws.Cells.Item(660, 15).Value = 20000
$ws.Cells.Item(660, 16).Value = 19750
$ws.Cells.Item(660, 19).Value = 988
$ws.Cells.Item(661, 4).Value = 44484
$ws.Cells.Item(661, 13).Value = 400
$ws.Cells.Item(661, 14).Value = 21000
$ws.Cells.Item(661, 15).Value = 22000
$ws.Cells.Item(661, 16).Value = 21500
$ws.Cells.Item(661, 19).Value = 1075
$ws.Cells.Item(662, 4).Value = 44665
$ws.Cells.Item(662, 14).Value = 17000
$ws.Cells.Item(662, 15).Value = 17000
$ws.Cells.Item(662, 16).Value = 17000
$ws.Cells.Item(662, 19).Value = 850
$ws.Cells.Item(663, 4).Value = 44665
$ws.Cells.Item(663, 13).Value = 160
$ws.Cells.Item(663, 14).Value = 18000
$ws.Cells.Item(663, 15).Value = 19000
$ws.Cells.Item(663, 16).Value = 18500
$ws.Cells.Item(663, 19).Value = 925
$ws.Cells.Item(664, 4).Value = 44452
$ws.Cells.Item(664, 13).Value = 80
$ws.Cells.Item(664, 14).Value = 20000
$ws.Cells.Item(664, 15).Value = 20000
$ws.Cells.Item(664, 16).Value = 20000
$ws.Cells.Item(664, 19).Value = 1000
$ws.Cells.Item(665, 4).Value = 44452
$ws.Cells.Item(665, 14).Value = 21000
$ws.Cells.Item(665, 15).Value = 22000
$ws.Cells.Item(665, 16).Value = 21500
$ws.Cells.Item(665, 19).Value = 1075
$ws.Cells.Item(666, 4).Value = 44510
$ws.Cells.Item(666, 13).Value = 160
$ws.Cells.Item(666, 14).Value = 16500
$ws.Cells.Item(666, 15).Value = 17000
$ws.Cells.Item(666, 16).Value = 16750
$ws.Cells.Item(666, 19).Value = 838
$ws.Cells.Item(667, 4).Value = 44510
$ws.Cells.Item(667, 13).Value = 240
$ws.Cells.Item(667, 14).Value = 18000
$ws.Cells.Item(667, 15).Value = 19000
$ws.Cells.Item(667, 16).Value = 18500
$ws.Cells.Item(667, 19).Value = 925
$ws.Cells.Item(668, 4).Value = 44189
$ws.Cells.Item(668, 13).Value = 680
$ws.Cells.Item(668, 15).Value = 10500
$ws.Cells.Item(668, 16).Value = 10279
$ws.Cells.Item(668, 19).Value = 514
$ws.Cells.Item(669, 4).Value = 44189
$ws.Cells.Item(669, 13).Value = 270
$ws.Cells.Item(669, 15).Value = 12500
$ws.Cells.Item(669, 16).Value = 12278
$ws.Cells.Item(669, 19).Value = 614
$ws.Cells.Item(670, 4).Value = 44701
$ws.Cells.Item(670, 13).Value = 120
$ws.Cells.Item(670, 14).Value = 10000
$ws.Cells.Item(670, 15).Value = 11000
$ws.Cells.Item(670, 16).Value = 10500
$ws.Cells.Item(670, 19).Value = 525
$ws.Cells.Item(671, 4).Value = 44701
$ws.Cells.Item(671, 13).Value = 160
$ws.Cells.Item(671, 14).Value = 12000
$ws.Cells.Item(671, 15).Value = 13000
$ws.Cells.Item(671, 16).Value = 12500
$ws.Cells.Item(671, 19).Value = 625
$ws.Cells.Item(672, 4).Value = 44516
$ws.Cells.Item(672, 13).Value = 180
$ws.Cells.Item(672, 14).Value = 16000
$ws.Cells.Item(672, 15).Value = 16000
$ws.Cells.Item(672, 16).Value = 16000
$ws.Cells.Item(672, 19).Value = 800
$ws.Cells.Item(673, 4).Value = 44516
$ws.Cells.Item(673, 13).Value = 240
$ws.Cells.Item(673, 14).Value = 17000
$ws.Cells.Item(673, 15).Value = 18000
$ws.Cells.Item(673, 16).Value = 17500
$ws.Cells.Item(673, 19).Value = 875
$ws.Cells.Item(674, 4).Value = 44186
$ws.Cells.Item(674, 13).Value = 740
$ws.Cells.Item(674, 14).Value = 10000
$ws.Cells.Item(674, 15).Value = 10500
$ws.Cells.Item(674, 16).Value = 10297
$ws.Cells.Item(674, 19).Value = 515
$ws.Cells.Item(675, 4).Value = 44186
$ws.Cells.Item(675, 13).Value = 350
$ws.Cells.Item(675, 14).Value = 11000
$ws.Cells.Item(675, 15).Value = 12000
$ws.Cells.Item(675, 16).Value = 11571
$ws.Cells.Item(675, 19).Value = 579
$ws.Cells.Item(676, 4).Value = 44463
$ws.Cells.Item(676, 13).Value = 150
$ws.Cells.Item(676, 14).Value = 14000
$ws.Cells.Item(676, 15).Value = 14000
$ws.Cells.Item(676, 16).Value = 14000
$ws.Cells.Item(676, 19).Value = 700
$ws.Cells.Item(677, 4).Value = 44463
$ws.Cells.Item(677, 13).Value = 400
$ws.Cells.Item(677, 14).Value = 15000
$ws.Cells.Item(677, 15).Value = 16000
$ws.Cells.Item(677, 16).Value = 15500
$ws.Cells.Item(677, 19).Value = 775

# Append two brand-new rows (678-679) duplicating the price records that
# used to be the last two rows (676-677) of the table, extending the sheet
# dimension from A1:T677 to A1:T679.

# Row 678
$ws.Cells.Item(678, 1).Value = 7
$ws.Cells.Item(678, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(678, 3).Value = "Ñuble"
$ws.Cells.Item(678, 5).Value = 16
$ws.Cells.Item(678, 6).Value = "Fruta"
$ws.Cells.Item(678, 7).Value = 100108
$ws.Cells.Item(678, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(678, 9).Value = 100108006
$ws.Cells.Item(678, 10).Value = "Plátano"
$ws.Cells.Item(678, 11).Value = "Sin especificar"
$ws.Cells.Item(678, 17).Value = "`$/caja 20 kilos"
$ws.Cells.Item(678, 18).Value = "Ecuador"
$ws.Cells.Item(678, 20).Value = 20
$ws.Cells.Item(678, 4).Value = 44382
$ws.Cells.Item(678, 4).NumberFormat = $ws.Cells.Item(676, 4).NumberFormat
$ws.Cells.Item(678, 12).Value = "Pintón"
$ws.Cells.Item(678, 13).Value = 80
$ws.Cells.Item(678, 14).Value = 10000
$ws.Cells.Item(678, 15).Value = 10000
$ws.Cells.Item(678, 16).Value = 10000
$ws.Cells.Item(678, 19).Value = 500

# Row 679
$ws.Cells.Item(679, 1).Value = 7
$ws.Cells.Item(679, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(679, 3).Value = "Ñuble"
$ws.Cells.Item(679, 5).Value = 16
$ws.Cells.Item(679, 6).Value = "Fruta"
$ws.Cells.Item(679, 7).Value = 100108
$ws.Cells.Item(679, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(679, 9).Value = 100108006
$ws.Cells.Item(679, 10).Value = "Plátano"
$ws.Cells.Item(679, 11).Value = "Sin especificar"
$ws.Cells.Item(679, 17).Value = "`$/caja 20 kilos"
$ws.Cells.Item(679, 18).Value = "Ecuador"
$ws.Cells.Item(679, 20).Value = 20
$ws.Cells.Item(679, 4).Value = 44382
$ws.Cells.Item(679, 4).NumberFormat = $ws.Cells.Item(676, 4).NumberFormat
$ws.Cells.Item(679, 12).Value = "Primera Pintón"
$ws.Cells.Item(679, 13).Value = 240
$ws.Cells.Item(679, 14).Value = 11000
$ws.Cells.Item(679, 15).Value = 12000
$ws.Cells.Item(679, 16).Value = 11500
$ws.Cells.Item(679, 19).Value = 575

